# Updated CVDs for the month
# Waldenburg Germany: set E4/E5 (cvd) to 0, and add a new "Commit/Forecast"
# data_source row (row 6) for Manufacturing Voluntary Turnover, mirroring
# the existing PY Actual / AOP rows for that metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Waldenburg Germany")

# cvd column (E) now has explicit 0 values on the existing
# "Manufacturing Voluntary Turnover" rows instead of being blank.
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0

# Duplicate row 5 into a new row 6 so the new row inherits the same
# number formatting/style as the rest of the table, then overwrite its
# contents with the new "Commit/Forecast" data_source values.
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "PES"
$ws.Range("B6").Value = "Commercial Systems"
$ws.Range("C6").Value = "Waldenburg Germany"
$ws.Range("D6").Value = "Manufacturing Voluntary Turnover"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "Commit/Forecast"
$ws.Range("G6:W6").Value = 0
